$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.387.68"
$ws.Range("E2").Value = "  -2.61%  "

$ws.Range("D3").Value = "2.428.93"
$ws.Range("E3").Value = "  -2.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "514.32"
$ws.Range("E5").Value = "  -3.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.53"
$ws.Range("E6").Value = "  -4.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("E8").Value = "  -2.93%  "

$ws.Range("D9").Value = "2.438.98"
$ws.Range("E9").Value = "  -2.69%  "

$ws.Range("E10").Value = "  -0.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0956"
$ws.Range("E11").Value = "  -5.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("E12").Value = "  -5.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.333"
$ws.Range("E13").Value = "  -3.31%  "

$ws.Range("D14").Value = "2.854.70"
$ws.Range("E14").Value = "  -3.22%  "

$ws.Range("D15").Value = "57.269.96"
$ws.Range("E15").Value = "  -2.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.63"
$ws.Range("E16").Value = "  -5.11%  "

$ws.Range("E17").Value = "  -4.24%  "

$ws.Range("D18").Value = "2.432.16"
$ws.Range("E18").Value = "  -2.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.44"
$ws.Range("E19").Value = "  -4.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "315.01"
$ws.Range("E20").Value = "  -2.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.10"
$ws.Range("E21").Value = "  -3.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.67"
$ws.Range("E23").Value = "  -3.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.66"
$ws.Range("E24").Value = "  -1.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.405"
$ws.Range("E25").Value = "  -3.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("E27").Value = "  -3.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.19"
$ws.Range("E28").Value = "  -4.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.18"
$ws.Range("E29").Value = "  +0.76%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.28"
$ws.Range("E30").Value = "  -2.82%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +3.40%  "

$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0717"
$ws.Range("E32").Value = "  -5.93%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.66"
$ws.Range("E33").Value = "  -4.71%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.72"
$ws.Range("E36").Value = "  -3.63%  "

$ws.Range("E37").Value = "  -6.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.91"
$ws.Range("E38").Value = "  -2.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.07"
$ws.Range("E39").Value = "  -2.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.45"
$ws.Range("E40").Value = "  -4.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.779"
$ws.Range("E41").Value = "  -2.70%  "

$ws.Range("E42").Value = "  -5.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "268.17"
$ws.Range("E43").Value = "  -4.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.94"
$ws.Range("E44").Value = "  -1.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.586"
$ws.Range("E45").Value = "  -2.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "121.01"
$ws.Range("E46").Value = "  -6.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0907"
$ws.Range("E47").Value = "  -1.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0483"
$ws.Range("E48").Value = "  -3.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.04"
$ws.Range("E49").Value = "  -4.74%  "

$ws.Range("E50").Value = "  -4.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.63"
$ws.Range("E51").Value = "  -3.87%  "
